$wb = $excel.ActiveWorkbook

# Sheet "展览" (overview/exhibition sheet) - rows 4-6 hold the F column counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 118
$wsExhibit.Range("F5").Value = 651
$wsExhibit.Range("F6").Value = 54

# Sheet "全部类型" (all types, combined listing) - same records appear on rows 5-7
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 118
$wsAll.Range("F6").Value = 651
$wsAll.Range("F7").Value = 54
